$d = $word.ActiveDocument

# The edit appends new material right after the very last paragraph of
# the document body ("The electrical energy released by FAOR ..."),
# and also gives that existing paragraph an explicit paragraph-mark
# run-properties block (<w:pPr><w:rPr>...</w:rPr></w:pPr>) that mirrors
# the formatting already used by its one run.
$last = $d.Paragraphs($d.Paragraphs.Count)
$r = $last.Range

# Preserve the existing paragraph's identity attributes (paraId / textId /
# rsid...) verbatim - we are only adding a <w:pPr> to it, not replacing it.
$lastParaXml = $r.WordOpenXML
$paraAttrs = ""
if ($lastParaXml -match '<w:p\s([^>]*)>') {
  $paraAttrs = " " + $matches[1]
}

# Reuse the formatting already on that paragraph's run for the new
# <w:pPr><w:rPr> block and for the first new (empty) paragraph.
$fontAscii = $r.Font.NameAscii
$fontCs = $r.Font.NameOther
$lang = $r.LanguageID

$wNs = "xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`""
$w14Ns = "xmlns:w14=`"http://schemas.microsoft.com/office/word/2010/wordml`""

$rPrMain = "<w:rPr><w:rFonts w:ascii=`"$fontAscii`" w:hAnsi=`"$fontAscii`" w:cs=`"$fontCs`"/><w:lang w:val=`"$lang`"/></w:rPr>"
$rPrEastAsia = "<w:rPr><w:rFonts w:ascii=`"$fontAscii`" w:eastAsiaTheme=`"minorEastAsia`" w:hAnsi=`"$fontAscii`" w:cs=`"$fontCs`"/><w:lang w:val=`"$lang`"/></w:rPr>"

$xml = "<w:p $wNs $w14Ns$paraAttrs><w:pPr>$rPrMain</w:pPr><w:r>$rPrMain<w:t>The electrical energy released by FAOR per electron can be estimated with the free-energy diagram of the direct pathway, as seen in figure x.</w:t></w:r></w:p>"
$xml += "<w:p $wNs><w:pPr>$rPrMain</w:pPr></w:p>"
$xml += "<w:p $wNs><w:pPr>$rPrEastAsia</w:pPr><w:r>$rPrEastAsia<w:t>FOAR highest U vs CHE: negative</w:t></w:r></w:p>"
$xml += "<w:p $wNs><w:pPr>$rPrEastAsia</w:pPr><w:r>$rPrEastAsia<w:t>OER highest U vs CHE: positive</w:t></w:r></w:p>"
$xml += "<w:p $wNs><w:pPr>$rPrEastAsia</w:pPr><w:r>$rPrEastAsia<w:t>The difference is the V_OC</w:t></w:r></w:p>"
$xml += "<w:p $wNs/>"

$r.InsertXML($xml)
Write-Output "Applied Garbage-Can fuel-cells edit."
